# Se agrega Script 0753 a la clase Tests_MiPortal
# Adds 13 new data rows (DEC_0758 .. DEC_0770) to the DataPool sheet,
# inserted above the trailing "TC / USUARIO / PASSWORD" template rows,
# which shifts those footer rows further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 13 blank rows at row 43 (pushes the old rows 43-47 down to 56-60),
# inheriting number formats/styles from the row above the insertion point
# (row 42), matching columns A-J: s=1,2,3,1,1,1,1,1,1,1.
$ws.Rows("43:55").Insert()

for ($i = 0; $i -lt 13; $i++) {
    $row = 43 + $i
    $num = 758 + $i
    $code = "DEC_0" + $num

    $ws.Cells.Item($row, 1).Value  = $code
    $ws.Cells.Item($row, 2).Value  = "18092588-0"
    # Leading apostrophe preserves the quote-prefix ("Text", quotePrefix)
    # style already present on column C (style index 3) instead of it
    # being reset to the plain text style when the value is assigned.
    $ws.Cells.Item($row, 3).Value  = "'sebA$1357"
    $ws.Cells.Item($row, 4).Value  = "SIN_DATO"
    $ws.Cells.Item($row, 5).Value  = "SIN_DATO"
    $ws.Cells.Item($row, 6).Value  = "SIN_DATO"
    $ws.Cells.Item($row, 7).Value  = "SIN_DATO"
    $ws.Cells.Item($row, 8).Value  = "SIN_DATO"
    $ws.Cells.Item($row, 9).Value  = "SIN_DATO"
    $ws.Cells.Item($row, 10).Value = "SIN_DATO"
}

# Update the saved view/selection to match the post-edit state.
[void]$ws.Range("C41").Select()
